$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Row 2 - NVDA
Set-TextCell "C2" '$103.83'
Set-TextCell "D2" '$3.83'
Set-TextCell "E2" '3.83%'
$ws.Range("F2").Value = 11
$ws.Range("G2").Value = 3
Set-TextCell "H2" '27.27%'
Set-TextCell "J2" '2025-08-08 08:36:13'

# Row 3 - AUR
$ws.Range("F3").Value = 7
$ws.Range("I3").Value = 2
Set-TextCell "J3" '2025-08-08 08:36:13'

# Row 4 - TSLA
Set-TextCell "J4" '2025-08-08 08:36:13'

# Row 5 - SOFI
$ws.Range("I5").Value = 3
Set-TextCell "J5" '2025-08-08 08:36:13'

# Row 6 - SOUN
Set-TextCell "C6" '$100.43'
Set-TextCell "D6" '$0.43'
Set-TextCell "E6" '0.43%'
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 1
Set-TextCell "H6" '12.50%'
Set-TextCell "J6" '2025-08-08 08:36:13'

# Row 7 - AMD
Set-TextCell "J7" '2025-08-08 08:36:13'

# Row 8 - AVGO
Set-TextCell "C8" '$110.91'
Set-TextCell "D8" '$10.91'
Set-TextCell "E8" '10.91%'
$ws.Range("F8").Value = 7
$ws.Range("G8").Value = 3
Set-TextCell "H8" '42.86%'
$ws.Range("I8").Value = 2
Set-TextCell "J8" '2025-08-08 08:36:13'

# Row 9 - CRCL
Set-TextCell "C9" '$140.89'
Set-TextCell "D9" '$40.89'
Set-TextCell "E9" '40.89%'
$ws.Range("F9").Value = 6
$ws.Range("G9").Value = 2
Set-TextCell "H9" '33.33%'
$ws.Range("I9").Value = 3
Set-TextCell "J9" '2025-08-08 08:36:13'

# Row 10 - BBAI
Set-TextCell "C10" '$100.53'
Set-TextCell "D10" '$0.53'
Set-TextCell "E10" '0.53%'
$ws.Range("F10").Value = 6
$ws.Range("G10").Value = 2
Set-TextCell "H10" '33.33%'
$ws.Range("I10").Value = 3
Set-TextCell "J10" '2025-08-08 08:36:13'

# Row 11 - SLDB
Set-TextCell "C11" '$101.27'
Set-TextCell "D11" '$1.27'
Set-TextCell "E11" '1.27%'
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 2
Set-TextCell "H11" '50.00%'
$ws.Range("I11").Value = 2
Set-TextCell "J11" '2025-08-08 08:36:13'
